$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.008, -0.448, -0.212, 97,   59,  -63),
    @(0.153,  -0.041,  0.388, -89,  55,    0),
    @(-0.297, -0.169,  0.243, -108, 25,   98),
    @(-0.173, -0.08,   0.599, -33,  21,  118),
    @(0.2,    -0.259,  0.555, 83,  -44,  103),
    @(0.156,  -0.002,  0.598, 58,   8,   165),
    @(-0.134, -0.322,  0.576, 31,  49,  -107),
    @(-0.061, -0.162,  0.519, 89, -16,   -10),
    @(0.205,  -0.351,  0.267, 97,  38,   129),
    @(-0.372,  0.004,  0.235, -40, 30,  -121)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
